$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.866.29'
$ws.Range("E2").Value = '  -1.13%  '

$ws.Range("D3").Value = '1.735.22'
$ws.Range("E3").Value = '  +0.97%  '

$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").Value = '311.09'
$ws.Range("E5").Value = '  -0.26%  '

$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.04%  '

$ws.Range("D7").Value = '0.5011'
$ws.Range("E7").Value = '  +9.48%  '

$ws.Range("D8").Value = '0.3559'
$ws.Range("E8").Value = '  +3.97%  '

$ws.Range("D9").Value = '42.13'
$ws.Range("E9").Value = '  +0.30%  '

$ws.Range("D10").Value = '0.07240'
$ws.Range("E10").Value = '  -0.16%  '

$ws.Range("D11").Value = '1.056'
$ws.Range("E11").Value = '  +1.31%  '

$ws.Range("D12").Value = '1.003'
$ws.Range("E12").Value = '  +0.06%  '

$ws.Range("D13").Value = '20.14'
$ws.Range("E13").Value = '  +2.03%  '

$ws.Range("D14").Value = '5.930'
$ws.Range("E14").Value = '  +1.87%  '

$ws.Range("D15").Value = '1.739.91'
$ws.Range("E15").Value = '  +0.94%  '

$ws.Range("D16").Value = '6.803'
$ws.Range("E16").Value = '  -0.76%  '

$ws.Range("D17").Value = '86.46'
$ws.Range("E17").Value = '  -2.31%  '

$ws.Range("D18").Value = '0.00001033'
$ws.Range("E18").Value = '  -0.57%  '

$ws.Range("D19").Value = '0.06429'
$ws.Range("E19").Value = '  +1.59%  '

$ws.Range("E20").Value = '  +0.01%  '

$ws.Range("D21").Value = '16.45'
$ws.Range("E21").Value = '  -0.39%  '

$ws.Range("D22").Value = '5.709'
$ws.Range("E22").Value = '  +1.72%  '

$ws.Range("D23").Value = '26.935.21'
$ws.Range("E23").Value = '  -1.05%  '

$ws.Range("D24").Value = '11.27'
$ws.Range("E24").Value = '  +3.99%  '

$ws.Range("D25").Value = '2.047'
$ws.Range("E25").Value = '  -3.85%  '

$ws.Range("D26").Value = '153.57'
$ws.Range("E26").Value = '  -0.76%  '

$ws.Range("D27").Value = '19.77'
$ws.Range("E27").Value = '  +2.74%  '

$ws.Range("D28").Value = '1.941.46'
$ws.Range("E28").Value = '  +1.11%  '

$ws.Range("D29").Value = '2.202'
$ws.Range("E29").Value = '  +3.31%  '

$ws.Range("D30").Value = '119.72'
$ws.Range("E30").Value = '  +0.27%  '

$ws.Range("D31").Value = '1.040'
$ws.Range("E31").Value = '  +1.52%  '

$ws.Range("D32").Value = '0.09507'
$ws.Range("E32").Value = '  +4.39%  '

$ws.Range("D33").Value = '3.586'
$ws.Range("E33").Value = '  -0.30%  '

$ws.Range("D34").Value = '5.338'
$ws.Range("E34").Value = '  +0.07%  '

$ws.Range("D35").Value = '0.02188'
$ws.Range("E35").Value = '  -0.12%  '

$ws.Range("D36").Value = '0.05812'
$ws.Range("E36").Value = '  -0.40%  '

$ws.Range("D37").Value = '11.02'
$ws.Range("E37").Value = '  -0.06%  '

$ws.Range("B38").Value = 'WEMIXTOKEN'
$ws.Range("C38").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D38").Value = '1.423'
$ws.Range("E38").Value = '  +1.37%  '

$ws.Range("B39").Value = 'Algorand'
$ws.Range("C39").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D39").Value = '0.1998'
$ws.Range("E39").Value = '  +0.18%  '

$ws.Range("D40").Value = '4.759'
$ws.Range("E40").Value = '  +0.18%  '

$ws.Range("D41").Value = '0.6028'
$ws.Range("E41").Value = '  +2.06%  '

$ws.Range("D42").Value = '1.107'
$ws.Range("E42").Value = '  -2.04%  '

$ws.Range("E43").Value = '  +1.71%  '

$ws.Range("D44").Value = '12.79'
$ws.Range("E44").Value = '  +0.57%  '

$ws.Range("D45").Value = '3.594'
$ws.Range("E45").Value = '  +0.21%  '

$ws.Range("D46").Value = '0.5636'
$ws.Range("E46").Value = '  +0.23%  '

$ws.Range("D47").Value = '119.88'
$ws.Range("E47").Value = '  +0.88%  '

$ws.Range("D48").Value = '1.845'
$ws.Range("E48").Value = '  -0.90%  '

$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '0.06661'
$ws.Range("E49").Value = '  +0.06%  '

$ws.Range("B50").Value = 'EOS'
$ws.Range("C50").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D50").Value = '1.095'
$ws.Range("E50").Value = '  +1.13%  '

$ws.Range("D51").Value = '1.002'
$ws.Range("E51").Value = '  +0.04%  '
